$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B104").Value = "ร้านขายของออนไลน์แบบไหน ที่ต้องจดทะเบียนพาณิชย์อิเล็กทรอนิกส์"
$ws.Range("B105").Value = "กฎหมายการขายของอิเล็กทรอนิกส์"
$ws.Range("B106").Value = "กฎหมายขายของอิเล็กทรอนิกส์"
$ws.Range("B107").Value = "จดทะเบียนอิเล็กทรอนิกส์"
$ws.Range("B108").Value = "อิเล็กทรอนิกส์"
$ws.Range("B110").Value = "ทะเบียนพาณิชย์อิเล็กทรอนิกส์"
$ws.Range("B111").Value = "ร้านแบบไหนต้องขดทะเบียนอิเล็กทรอนิกส์"
$ws.Range("B112").Value = "ร้านขายของออนไลน์แบบไหน "
$ws.Range("B113").Value = "ร้านที่ต้องจดทะเบียนอิเล็กทรอนิกส์"
$ws.Range("B114").Value = "ทะเบียนพาณิชย์อิเล็กทรอนิกส์"
$ws.Range("B115").Value = "จดทะเบียนพาณิชย์อิเล็กทรอนิกส์"
$ws.Range("B116").Value = "การจดทะเบียนพาณิชย์อิเล็กทรอนิกส์"
$ws.Range("B117").Value = "การจดทะเบียนอิเล็กทรอนิกส์"
$ws.Range("B118").Value = "DBD Registered"
$ws.Range("B119").Value = "dbd registered"
